{"js": "// Replacements: [findText, replaceText]\nconst replacements = [\n  [\"2024-05-03 Friday\", \"2024-05-04 Saturday\"],\n  [\"105\\u00F79=\", \"613\\u00F78=\"],\n  [\"847\\u00F79=\", \"358\\u00F78=\"],\n  [\"900\\u00F77=\", \"847\\u00F73=\"],\n  [\"954\\u00F78=\", \"749\\u00F77=\"],\n  [\"514\\u00F79=\", \"765\\u00F79=\"],\n  [\"733\\u00F77=\", \"343\\u00F76=\"],\n  [\"871\\u00F74=\", \"990\\u00F72=\"],\n  [\"427\\u00F73=\", \"836\\u00F72=\"],\n  [\"504\\u00F76=\", \"293\\u00F75=\"],\n  [\"427\\u00F72=\", \"903\\u00F75=\"],\n  [\"627\\u00F72=\", \"301\\u00F78=\"],\n  [\"259\\u00F77=\", \"545\\u00F75=\"],\n  [\"706\\u00F74=\", \"250\\u00F75=\"],\n  [\"471\\u00F73=\", \"390\\u00F78=\"],\n  [\"209\\u00F78=\", \"850\\u00F75=\"],\n  [\"222\\u00F76=\", \"883\\u00F78=\"],\n  [\"905\\u00F73=\", \"960\\u00F79=\"],\n  [\"505\\u00F78=\", \"575\\u00F78=\"],\n  [\"909\\u00F77=\", \"548\\u00F76=\"],\n  [\"476\\u00F76=\", \"151\\u00F76=\"],\n  [\"620\\u00F76=\", \"114\\u00F74=\"],\n  [\"253\\u00F79=\", \"651\\u00F72=\"],\n  [\"367\\u00F77=\", \"510\\u00F73=\"],\n  [\"829\\u00F72=\", \"153\\u00F74=\"],\n  [\"873\\u00F74=\", \"162\\u00F75=\"],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = context.document.body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-05-03 Friday\", \"2024-05-04 Saturday\"),\n    @(\"105\u00f79=\", \"613\u00f78=\"),\n    @(\"847\u00f79=\", \"358\u00f78=\"),\n    @(\"900\u00f77=\", \"847\u00f73=\"),\n    @(\"954\u00f78=\", \"749\u00f77=\"),\n    @(\"514\u00f79=\", \"765\u00f79=\"),\n    @(\"733\u00f77=\", \"343\u00f76=\"),\n    @(\"871\u00f74=\", \"990\u00f72=\"),\n    @(\"427\u00f73=\", \"836\u00f72=\"),\n    @(\"504\u00f76=\", \"293\u00f75=\"),\n    @(\"427\u00f72=\", \"903\u00f75=\"),\n    @(\"627\u00f72=\", \"301\u00f78=\"),\n    @(\"259\u00f77=\", \"545\u00f75=\"),\n    @(\"706\u00f74=\", \"250\u00f75=\"),\n    @(\"471\u00f73=\", \"390\u00f78=\"),\n    @(\"209\u00f78=\", \"850\u00f75=\"),\n    @(\"222\u00f76=\", \"883\u00f78=\"),\n    @(\"905\u00f73=\", \"960\u00f79=\"),\n    @(\"505\u00f78=\", \"575\u00f78=\"),\n    @(\"909\u00f77=\", \"548\u00f76=\"),\n    @(\"476\u00f76=\", \"151\u00f76=\"),\n    @(\"620\u00f76=\", \"114\u00f74=\"),\n    @(\"253\u00f79=\", \"651\u00f72=\"),\n    @(\"367\u00f77=\", \"510\u00f73=\"),\n    @(\"829\u00f72=\", \"153\u00f74=\"),\n    @(\"873\u00f74=\", \"162\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $find\n    $range.Find.Replacement.Text = $replace\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1  # wdFindContinue\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute([ref]$find, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$replace, [ref]2)\n}\n"}
